# Scheduled runner update: refresh Universalis market-price snapshots
# (currentAveragePrice / NQ / HQ, Leve prices, and derived profit columns)
# across the Golem_Profits Leve-crafting workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19: Unbreak My Heart / Roof Tile
$ws.Range("H19").Value = 267.66666
$ws.Range("I19").Value = 328.44446
$ws.Range("J19").Value = 85.333336
$ws.Range("K19").Value = 328.44446
$ws.Range("L19").Value = 85.333336
$ws.Range("M19").Value = -153.44446
$ws.Range("N19").Value = -435.333336

# Row 64: Forged from the Void / Void Glue
$ws.Range("H64").Value = 3666.3333
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3666.3333
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 3666.3333
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -4162.3333

# Row 67: Dodging the Draft (L) / Void Glue
$ws.Range("H67").Value = 3666.3333
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3666.3333
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 3666.3333
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -5382.3333

# Row 69: Steeling the Knife, Steeling the Mind / Grade 1 Mind Dissolvent
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").ClearContents()

# Row 72: Surgical Substitution (L) / Grade 1 Mind Dissolvent
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").ClearContents()

# Row 74: Adhesive of Antipathy / Wing Glue
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()

# Row 77: It's Gonna Grow Back (L) / Wing Glue
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()

# Row 107: Another Man's Ink / Enchanted Truegold Ink
$ws.Range("H107").Value = 397.64285
$ws.Range("I107").Value = 388
$ws.Range("K107").Value = 388
$ws.Range("M107").Value = 1532

$ws = $wb.Worksheets.Item("ARM")
# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 1865.25
$ws.Range("I45").Value = 2075.6667
$ws.Range("K45").Value = 2075.6667
$ws.Range("M45").Value = -1698.6667

# Row 63: Rivets Run through It / Mythrite Rivets
$ws.Range("H63").Value = 1862.3334
$ws.Range("I63").Value = 1833.6
$ws.Range("K63").Value = 1833.6
$ws.Range("M63").Value = -1147.6

# Row 66: A Riveting Revival (L) / Mythrite Rivets
$ws.Range("H66").Value = 1862.3334
$ws.Range("I66").Value = 1833.6
$ws.Range("K66").Value = 9168
$ws.Range("M66").Value = -5736

# Row 139: Backing up My Words / Titanium Gold Thornplate of Fending
$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280

$ws = $wb.Worksheets.Item("BSM")
# Row 29: Powderpost Derby / Initiate's Saw
$ws.Range("H29").Value = 1224.2667
$ws.Range("I29").Value = 1309.1428
$ws.Range("J29").Value = 1150
$ws.Range("K29").Value = 1309.1428
$ws.Range("L29").Value = 1150
$ws.Range("M29").Value = -1020.1428
$ws.Range("N29").Value = -1728

# Row 30: The Devil's Workshop / Brass Viking Sword
$ws.Range("H30").Value = 50000
$ws.Range("J30").Value = 50000
$ws.Range("L30").Value = 50000
$ws.Range("N30").Value = -50250

$ws = $wb.Worksheets.Item("CRP")
# Row 7: Gridania's Got Talent / Maple Lumber
$ws.Range("H7").Value = 130.41667
$ws.Range("I7").Value = 177.375
$ws.Range("K7").Value = 177.375
$ws.Range("M7").Value = -64.375

# Row 35: Storm of Swords / Elm Macuahuitl
$ws.Range("H35").Value = 6516
$ws.Range("I35").Value = 3145
$ws.Range("K35").Value = 3145
$ws.Range("M35").Value = -2851

# Row 80: The Long Armillae of the Law / Hallowed Chestnut Armillae
$ws.Range("H80").Value = 34999.5
$ws.Range("J80").Value = 34999.5
$ws.Range("L80").Value = 34999.5
$ws.Range("N80").Value = -37245.5

# Row 83: Wooden Ambitions (L) / Hallowed Chestnut Armillae
$ws.Range("H83").Value = 34999.5
$ws.Range("J83").Value = 34999.5
$ws.Range("L83").Value = 104998.5
$ws.Range("N83").Value = -116230.5

# Row 97: Wood That You Could / Larch Bracelets
$ws.Range("H97").Value = 61666.668
$ws.Range("J97").Value = 61666.668
$ws.Range("L97").Value = 61666.668
$ws.Range("N97").Value = -63648.668

$ws = $wb.Worksheets.Item("CUL")
# Row 141: Ocean Explosion / Acqua Pazza
$ws.Range("H141").Value = 1797.3334
$ws.Range("I141").Value = 1557.8
$ws.Range("J141").Value = 2995
$ws.Range("K141").Value = 4673.4
$ws.Range("L141").Value = 8985
$ws.Range("M141").Value = 506.6000000000004
$ws.Range("N141").Value = -19345

$ws = $wb.Worksheets.Item("GSM")
# Row 7: Water of Life / Copper Rings
$ws.Range("H7").Value = 252.75
$ws.Range("I7").Value = 3.142857
$ws.Range("K7").Value = 3.142857
$ws.Range("M7").Value = 108.857143

# Row 8: Gods of Small Things / Copper Earrings
$ws.Range("H8").Value = 252.75
$ws.Range("I8").Value = 3.142857
$ws.Range("K8").Value = 3.142857
$ws.Range("M8").Value = 135.857143

# Row 26: Perk of Fiction / Coral Ring
$ws.Range("H26").Value = 8400
$ws.Range("I26").Value = 8400
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 8400
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -8120
$ws.Range("N26").ClearContents()

# Row 50: Coral on My Mind / Red Coral Ring
$ws.Range("H50").Value = 8400
$ws.Range("I50").Value = 8400
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 8400
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -7902
$ws.Range("N50").ClearContents()

# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# Row 106: Choker in the Clutch / Palladium Choker of Aiming
$ws.Range("H106").Value = 29994
$ws.Range("J106").Value = 29994
$ws.Range("L106").Value = 29994
$ws.Range("N106").Value = -32518

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban / Leather
$ws.Range("H7").Value = 8397.366
$ws.Range("I7").Value = 8590.076999999999
$ws.Range("J7").Value = 7144.75
$ws.Range("K7").Value = 8590.076999999999
$ws.Range("L7").Value = 7144.75
$ws.Range("M7").Value = -8478.076999999999
$ws.Range("N7").Value = -7368.75

# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 1738.7
$ws.Range("I22").Value = 1296
$ws.Range("J22").Value = 1928.4286
$ws.Range("K22").Value = 1296
$ws.Range("L22").Value = 1928.4286
$ws.Range("M22").Value = -1001
$ws.Range("N22").Value = -2518.4286

# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 1738.7
$ws.Range("I27").Value = 1296
$ws.Range("J27").Value = 1928.4286
$ws.Range("K27").Value = 1296
$ws.Range("L27").Value = 1928.4286
$ws.Range("M27").Value = -1189
$ws.Range("N27").Value = -2142.4286

# Row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 4833.3335
$ws.Range("I68").Value = 6000
$ws.Range("K68").Value = 6000
$ws.Range("M68").Value = -5251

# Row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 4833.3335
$ws.Range("I71").Value = 6000
$ws.Range("K71").Value = 30000
$ws.Range("M71").Value = -26256

# Row 82: Trainin' the Neck / Dragon Leather
$ws.Range("H82").Value = 3078.2727
$ws.Range("I82").Value = 2123
$ws.Range("J82").Value = 4750
$ws.Range("K82").Value = 2123
$ws.Range("L82").Value = 4750
$ws.Range("M82").Value = -1762
$ws.Range("N82").Value = -5472

# Row 85: Training Is Only Skintight (L) / Dragon Leather
$ws.Range("H85").Value = 3078.2727
$ws.Range("I85").Value = 2123
$ws.Range("J85").Value = 4750
$ws.Range("K85").Value = 2123
$ws.Range("L85").Value = 4750
$ws.Range("M85").Value = -875
$ws.Range("N85").Value = -7246

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 3686.0667
$ws.Range("I122").Value = 3327.2856
$ws.Range("K122").Value = 9981.856800000001
$ws.Range("M122").Value = -7531.856800000001

# Row 126: Battered Books / Saiga Leather
$ws.Range("H126").Value = 8397.366
$ws.Range("I126").Value = 8590.076999999999
$ws.Range("J126").Value = 7144.75
$ws.Range("K126").Value = 25770.231
$ws.Range("L126").Value = 21434.25
$ws.Range("M126").Value = -23300.231
$ws.Range("N126").Value = -26374.25

# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 7666.6665
$ws.Range("I132").Value = 7666.6665
$ws.Range("K132").Value = 22999.9995
$ws.Range("M132").Value = -20469.9995

$ws = $wb.Worksheets.Item("WVR")
# Row 51: After the Smock-down / Linen Smock
$ws.Range("H51").Value = 20038
$ws.Range("J51").Value = 20038
$ws.Range("L51").Value = 20038
$ws.Range("N51").Value = -21058

# Row 95: Duress Rehearsal / Ruby Cotton Fingerless Gloves of Casting
$ws.Range("H95").Value = 17250
$ws.Range("J95").Value = 17250
$ws.Range("L95").Value = 17250
$ws.Range("N95").Value = -22742

# Row 102: Don't Sweat the Role / Serge Turban of Crafting
$ws.Range("H102").Value = 130000
$ws.Range("J102").Value = 130000
$ws.Range("L102").Value = 130000
$ws.Range("N102").Value = -136490

# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 6000
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -10900

# Row 126: A Polished Purchase / Snow Linen
$ws.Range("H126").Value = 5603
$ws.Range("I126").Value = 3306
$ws.Range("K126").Value = 9918
$ws.Range("M126").Value = -7448

# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 8471.666999999999
$ws.Range("I132").Value = 7189
$ws.Range("K132").Value = 21567
$ws.Range("M132").Value = -19037
